$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (Shrikt / 405,HR,Hall)
$ws.Range("A4").Value = "Shrikt"
$ws.Range("B4").Value = "405,HR,Hall"

# Update the active cell selection to match the saved view state
$ws.Range("D9").Select()
